# Preenche a planilha com os dados formatados, incluindo a função de formatação do CPF
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sobrenome (coluna B) para todas as 11 linhas de dados
for ($n = 1; $n -le 11; $n++) {
    $ws.Cells.Item($n, 2).Value2 = "Ribeiro"
}

# Primeiro CPF (coluna C, linha 1)
$ws.Cells.Item(1, 3).Value2 = "115.853.176-16"

# Nomes (coluna A) para todas as 11 linhas
for ($n = 1; $n -le 11; $n++) {
    $ws.Cells.Item($n, 1).Value2 = "Alexandre$n"
}

# CPFs formatados restantes (coluna C, linhas 2 a 11)
for ($n = 2; $n -le 11; $n++) {
    $cpf = 15 + $n
    $ws.Cells.Item($n, 3).Value2 = "115.853.176-$cpf"
}

# Remove a antiga linha 12 (cabecalho duplicado/estilizado)
$ws.Rows.Item(12).Delete()

# Sublinha a celula B11
$ws.Range("B11").Font.Underline = 2

# Celula F13 vazia, com a mesma formatacao (sublinhado)
$ws.Range("F13").Font.Underline = 2

# Larguras de coluna
$ws.Columns.Item(1).ColumnWidth = 15.67
$ws.Columns.Item(2).ColumnWidth = 6.67

# Configuracao de pagina
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selecao final
$ws.Range("F13").Select() | Out-Null
